$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'59.335.68"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.58%  "
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'2.527.54"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.28%  "
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'535.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.04%  "
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'140.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -2.52%  "
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("E7").Value = "'  +0.32%  "
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("E8").Value = "'  -2.00%  "
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").Value = "'2.530.62"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +0.08%  "
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("E10").Value = "'  -0.35%  "
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("E11").Value = "'  +1.70%  "
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").Value = "'5.38"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -1.18%  "
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("E13").Value = "'  +0.63%  "
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").Value = "'2.975.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +1.30%  "
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("E15").Value = "'  -2.72%  "
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'59.269.49"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +0.67%  "
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").Value = "'0.0000141"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +1.13%  "
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").Value = "'2.553.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +1.58%  "
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("D19").Value = "'10.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -2.65%  "
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("E20").Value = "'  -1.29%  "
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").Value = "'321.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -0.21%  "
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("E22").Value = "'  +0.25%  "
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").Value = "'5.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +1.41%  "
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").Value = "'62.44"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +1.09%  "
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("E25").Value = "'  -4.15%  "
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").Value = "'0.165"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +2.25%  "
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("E27").Value = "'  +0.52%  "
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").Value = "'7.84"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +1.18%  "
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("E29").Value = "'  -0.48%  "
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("E30").Value = "'  -0.84%  "
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("E31").Value = "'  +0.38%  "
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("D32").Value = "'161.12"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.80%  "
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("E33").Value = "'  +0.27%  "
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").Value = "'1.14"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -4.23%  "
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("E35").Value = "'  -0.68%  "
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("D36").Value = "'18.51"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -0.15%  "
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("E37").Value = "'  -2.96%  "
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("E38").Value = "'  -2.45%  "
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").Value = "'37.06"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +0.55%  "
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("E40").Value = "'  -0.66%  "
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("B41").Value = "'Bittensor"
$ws.Range("B41").ClearFormats()
$ws.Range("C41").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C41").ClearFormats()
$ws.Range("D41").Value = "'286.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -5.64%  "
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("B42").Value = "'SuiNetwork"
$ws.Range("B42").ClearFormats()
$ws.Range("C42").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C42").ClearFormats()
$ws.Range("D42").Value = "'0.804"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -2.06%  "
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("B43").Value = "'RenderToken"
$ws.Range("B43").ClearFormats()
$ws.Range("C43").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C43").ClearFormats()
$ws.Range("D43").Value = "'5.29"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -5.33%  "
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").Value = "'0.997"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.21%  "
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").Value = "'0.601"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +0.47%  "
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("E46").Value = "'  +0.59%  "
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47").Value = "'124.17"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.89%  "
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("E48").Value = "'  -0.24%  "
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("E49").Value = "'  +0.07%  "
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("E50").Value = "'  -1.62%  "
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("E51").Value = "'  -2.20%  "
$ws.Range("E51").ClearFormats()

